# Updated cryptos list on Wed Jul 24 14:47:00 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto listing, and fixes the swapped Fetch.AI/Aptos rows (34 & 35),
# which also picked up new Price/Volume values.
#
# Note: several Price values look like plain numbers (e.g. "13.90",
# "0.0700", "1.00"); Excel would silently coerce those to numeric cells
# and strip the formatting-significant trailing zeros, so NumberFormat is
# forced to Text ("@") right before such assignments to keep them as the
# original string values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.945.10"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "3.443.41"
$ws.Range("E3").Value = "  -1.43%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.80"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.86"
$ws.Range("E6").Value = "  +1.41%  "

$ws.Range("E7").Value = "  +4.87%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "3.442.69"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  +0.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").Value = "4.046.53"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.09"
$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").Value = "66.908.44"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "3.443.71"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.97"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.42"
$ws.Range("E21").Value = "  -2.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.68"
$ws.Range("E22").Value = "  -2.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.59"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("E24").Value = "  +7.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.537"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("E28").Value = "  +1.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.89"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.71"
$ws.Range("E32").Value = "  -3.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.01%  "

# Rows 34 & 35: Fetch.AI and Aptos swapped places, with refreshed data.
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.12"
$ws.Range("E34").Value = "  -1.42%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.29"
$ws.Range("E35").Value = "  -3.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.58"
$ws.Range("E36").Value = "  -0.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.06"
$ws.Range("E37").Value = "  +1.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.882"
$ws.Range("E38").Value = "  -1.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.89"
$ws.Range("E39").Value = "  -6.95%  "

$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.51"
$ws.Range("E42").Value = "  -0.73%  "

$ws.Range("D43").Value = "2.752.58"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.38"
$ws.Range("E44").Value = "  -1.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0700"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.69"
$ws.Range("E46").Value = "  +3.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "339.51"
$ws.Range("E47").Value = "  +6.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.24"
$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0288"
$ws.Range("E49").Value = "  -3.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.105"
$ws.Range("E50").Value = "  +2.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.06"
$ws.Range("E51").Value = "  +1.97%  "
